$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price strings in column D stay as text (matches original inlineStr formatting)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.180.53"
$ws.Range("E2").Value = "  -2.13%  "
$ws.Range("D3").Value = "1.578.46"
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  -0.54%  "
$ws.Range("D5").Value = "208.85"
$ws.Range("E5").Value = "  -1.43%  "
$ws.Range("E6").Value = "  -2.99%  "
$ws.Range("E7").Value = "  -0.52%  "
$ws.Range("E8").Value = "  -1.52%  "
$ws.Range("E9").Value = "  -1.10%  "
$ws.Range("D10").Value = "19.50"
$ws.Range("E10").Value = "  -0.62%  "
$ws.Range("D11").Value = "0.0843"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("D12").Value = "1.799.14"
$ws.Range("D13").Value = "1.598.06"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("D14").Value = "4.06"
$ws.Range("E14").Value = "  -0.38%  "
$ws.Range("E15").Value = "  -2.04%  "
$ws.Range("E16").Value = "  -1.21%  "
$ws.Range("D17").Value = "26.179.34"
$ws.Range("E17").Value = "  -2.08%  "
$ws.Range("E18").Value = "  -1.96%  "
$ws.Range("D19").Value = "7.37"
$ws.Range("E19").Value = "  +2.55%  "
$ws.Range("D20").Value = "209.09"
$ws.Range("E20").Value = "  -0.09%  "
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "4.25"
$ws.Range("E22").Value = "  -1.51%  "
$ws.Range("E23").Value = "  -2.94%  "
$ws.Range("D24").Value = "8.82"
$ws.Range("E24").Value = "  -2.57%  "
$ws.Range("D25").Value = "144.13"
$ws.Range("E25").Value = "  +0.40%  "
$ws.Range("E26").Value = "  -0.51%  "
$ws.Range("D27").Value = "7.00"
$ws.Range("E27").Value = "  -1.92%  "
$ws.Range("E28").Value = "  -1.88%  "
$ws.Range("D29").Value = "15.23"
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("E30").Value = "  -0.29%  "
$ws.Range("E31").Value = "  -1.33%  "
$ws.Range("D32").Value = "3.22"
$ws.Range("E32").Value = "  -1.73%  "
$ws.Range("D33").Value = "2.99"
$ws.Range("E33").Value = "  +0.81%  "
$ws.Range("D34").Value = "1.283.35"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("B36").Value = "ImmutableX"
$ws.Range("C36").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D36").Value = "0.606"
$ws.Range("E36").Value = "  +2.98%  "
$ws.Range("E37").Value = "  -1.65%  "
$ws.Range("B38").Value = "WEMIXToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D38").Value = "1.11"
$ws.Range("E38").Value = "  -9.24%  "
$ws.Range("B39").Value = "VeChain"
$ws.Range("C39").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D39").Value = "0.0166"
$ws.Range("E39").Value = "  -1.95%  "
$ws.Range("E40").Value = "  -2.00%  "
$ws.Range("E41").Value = "  -0.51%  "
$ws.Range("D42").Value = "5.59"
$ws.Range("E42").Value = "  +2.07%  "
$ws.Range("D43").Value = "0.764"
$ws.Range("E43").Value = "  -1.78%  "
$ws.Range("E44").Value = "  -3.14%  "
$ws.Range("D45").Value = "62.23"
$ws.Range("E45").Value = "  -0.75%  "
$ws.Range("D46").Value = "1.712.99"
$ws.Range("E46").Value = "  -1.62%  "
$ws.Range("D47").Value = "88.50"
$ws.Range("E47").Value = "  -2.09%  "
$ws.Range("E48").Value = "  -0.28%  "
$ws.Range("E49").Value = "  -4.18%  "
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("E51").Value = "  -1.43%  "
